# Insert a new data row at row 30 (pushes existing rows 30-60 down to 31-61)
# and populate it with the new record's values, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 30, shifting rows 30..60 down to 31..61.
$ws.Rows(30).Insert()

# Populate the newly inserted row 30 with the new record.
$ws.Cells.Item(30, 1).Value2  = 5
$ws.Cells.Item(30, 2).Value2  = "Macroferia Regional de Talca"
$ws.Cells.Item(30, 3).Value2  = "Maule"
$ws.Cells.Item(30, 4).Value2  = 44512
$ws.Cells.Item(30, 5).Value2  = 7
$ws.Cells.Item(30, 6).Value2  = 100112022
$ws.Cells.Item(30, 7).Value2  = "Arveja Verde"
$ws.Cells.Item(30, 8).Value2  = "Sin especificar"
$ws.Cells.Item(30, 9).Value2  = "Primera"
$ws.Cells.Item(30, 10).Value2 = 500
$ws.Cells.Item(30, 11).Value2 = 15000
$ws.Cells.Item(30, 12).Value2 = 15000
$ws.Cells.Item(30, 13).Value2 = 15000
$ws.Cells.Item(30, 14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(30, 15).Value2 = "Región del Maule"
$ws.Cells.Item(30, 16).Value2 = 600
$ws.Cells.Item(30, 17).Value2 = 25
$ws.Cells.Item(30, 18).Value2 = "Hortaliza"

# Apply the same date-cell style/format used by the other rows' "Fecha" column (D).
$ws.Cells.Item(30, 4).NumberFormat = $ws.Cells.Item(31, 4).NumberFormat
